$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column A to fit the longer indicator names (stored width "42")
$ws.Columns.Item(1).ColumnWidth = 41.166666666666664

# Populate the new indicator-metadata rows (3-10)
# Row 3
$ws.Range('A3').Value = 'BSB_Shelf_Water_Volume_South'
$ws.Range('B3').Value = 'data/BSB_Shelf_Water_Volume_South.txt'
$ws.Range('C3').Value = 'Shelf water volume (km^3) for the Southern stock subunit. Shelf water volume is a measure of the volume of water inshore of the shelf-slope front, a narrow transition region between masses of cool, low salinity shelf water and warm, high salinity slope water. This in-situ data was derived from CTD data from NEFSC surveys and is maintained by Paula Fratantoni.'
$ws.Range('D3').Value = 'There has been no winter sampling since 2021, no trend for 2024.'
$ws.Range('E3').Value = 'There has been no winter sampling since 2021, no trend for 2024.'
$ws.Range('F3').Value = 'Shelf water volume is a proxy for suitable winter habitat; higher shelf water volume indicates less suitable habitat, potentially leading to northern fish migrating into the southern subregion.'
$ws.Range('G3').Value = 'Fratantoni PS, Holzworth-Davis T, Taylor MH. 2015. Description of oceanographic conditions on the Northeast US Continental Shelf during 2014. US Dept Commer, Northeast Fisheries Science Center. Ref Doc. 15-21; 41 p. Available at: http://www.nefsc.noaa.gov/publications/doi:10.7289/V5SQ8XD2'
$ws.Range('H3').Value = 'SOUTH'

# Row 4
$ws.Range('A4').Value = 'BSB_Shelf_Water_Volume_North'
$ws.Range('B4').Value = 'data/BSB_Shelf_Water_Volume_North.txt'
$ws.Range('C4').Value = 'Shelf water volume (km^3) for the NoBSB_Shelf_Water_Volume_Northrthern stock subunit. Shelf water volume is a measure of the volume of water inshore of the shelf-slope front, a narrow transition region between masses of cool, low salinity shelf water and warm, high salinity slope water. This in-situ data was derived from CTD data from NEFSC surveys and is maintained by Paula Fratantoni.'
$ws.Range('D4').Value = 'There has been no winter sampling since 2021, no trend for 2024.'
$ws.Range('E4').Value = 'There has been no winter sampling since 2021, no trend for 2024.'
$ws.Range('F4').Value = 'Shelf water volume is a proxy for suitable winter habitat; higher shelf water volume indicates less suitable habitat, potentially leading to northern fish migrating into the southern subregion.'
$ws.Range('G4').Value = 'Fratantoni PS, Holzworth-Davis T, Taylor MH. 2015. Description of oceanographic conditions on the Northeast US Continental Shelf during 2014. US Dept Commer, Northeast Fisheries Science Center. Ref Doc. 15-21; 41 p. Available at: http://www.nefsc.noaa.gov/publications/doi:10.7289/V5SQ8XD2'
$ws.Range('H4').Value = 'NORTH'

# Row 5
$ws.Range('A5').Value = 'BSB_mrip_trips'
$ws.Range('B5').Value = 'data/bsb_rec_trips.txt'
$ws.Range('C5').Value = 'Total annual recreational black sea bass fishing trips for both North and South subunits. Data from NOAA Fisheries’ Marine Recreational Information Program (MRIP).'
$ws.Range('D5').Value = 'Recent trip numbers are near an all-time high, but may have decreased from 2023.'
$ws.Range('E5').Value = 'Catch generally reflects trip patterns. High regulatory complexity is likely contributing to recreational fishing trends.'
$ws.Range('F5').Value = 'Black sea bass is an important Mid-Atlantic stock with high recreational engagement.'
$ws.Range('H5').Value = 'ALL'

# Row 6
$ws.Range('A6').Value = 'BSB_mrip_landings'
$ws.Range('B6').Value = 'data/bsb_rec_landings.txt'
$ws.Range('C6').Value = 'Total annual recreational landings of black sea bass for both North and South subunits. Data from NOAA Fisheries’ Marine Recreational Information Program (MRIP)'
$ws.Range('D6').Value = 'Recreational landings have decreased from 2023 but remain near the long-term average.'
$ws.Range('E6').Value = 'The recreational black sea bass fishery has a catch-and-release component, and management measures are being implemented to reduce recreational harvest.'
$ws.Range('F6').Value = 'Black sea bass is an important Mid-Atlantic stock with high recreational engagement.'
$ws.Range('H6').Value = 'ALL'

# Row 7
$ws.Range('A7').Value = 'BSB_Commercial_Revenue'
$ws.Range('B7').Value = 'data/bsb_com_revenue.txt'
$ws.Range('C7').Value = 'Black sea bass commercial revenue (2024 USD)'
$ws.Range('D7').Value = 'Commercial revenue for black sea bass has increased from 2023 and remains well above the long term average'
$ws.Range('E7').Value = 'Black sea bass has high commercial value that continues to increase in 2024 despite fewer active vessels in the fleet.'
$ws.Range('F7').Value = 'Commercial revenue per vessel has an overall increasing trend, despite decreases in both total landings and average price ($/lb).'
$ws.Range('H7').Value = 'ALL'

# Row 8
$ws.Range('A8').Value = 'BSB_Commercial_Revenue_Per_Vessel'
$ws.Range('B8').Value = 'data/bsb_com_revenue_per_vessel.txt'
$ws.Range('C8').Value = 'Black sea bass commercial revenue per vessel (2024 USD)'
$ws.Range('D8').Value = 'Commercial revenue for black sea bass has increased from 2023 and remains well above the long term average'
$ws.Range('E8').Value = 'Black sea bass has high commercial value that continues to increase in 2024 despite fewer active vessels in the fleet.'
$ws.Range('F8').Value = 'Commercial revenue per vessel has an overall increasing trend, despite decreases in both total landings and average price ($/lb).'
$ws.Range('H8').Value = 'ALL'

# Row 9
$ws.Range('A9').Value = 'BSB_Commercial_Vessels'
$ws.Range('B9').Value = 'data/bsb_com_vessels.txt'''
$ws.Range('C9').Value = 'Number of commercial fishing vessels targeting black sea bass'
$ws.Range('D9').Value = 'Number of commercial vessels has decreased slightly from 2023 and remains well below the long term average.'
$ws.Range('E9').Value = 'A decrease in targeted black sea bass trips coincides with decreased catch and landings in 2024.'
$ws.Range('F9').Value = 'The number of active vessels has been decreasing since 2017, which could impact revenue distributions and fleet composition."'
$ws.Range('H9').Value = 'ALL'

# Row 10
$ws.Range('A10').Value = 'BSB_Winter_Bottom_Temperature_South'
$ws.Range('B10').Value = 'data/bsb_winter_bottom_temperature_south.txt'
$ws.Range('C10').Value = 'Winter (Feb-Mar) bottom temperature in the black sea bass South stock region. Hubert''s data product is a composite before 1993, and from 1993-2019 it is the same as GLORYS. 2020-2024 data are pulled directly from GLORYS. The GLORYS12V1 product is the CMEMS global ocean eddy-resolving (1/12? horizontal resolution, 50 vertical levels) reanalysis.'
$ws.Range('D10').Value = 'Bottom temperatures in 2024 are decreasing relative to recent years, but still within 1 sd of the mean.'
$ws.Range('E10').Value = 'Cold winter temperatures in the Northwest Atlantic (north of Hudson Canyon).'
$ws.Range('F10').Value = 'Cold winter temperatures may increase the mortality of young-of-the-year fish, resulting in smaller year classes. Additionally, cold temperatures can cause northern fish to move into the southern subregion, leading to potential misallocation of catch between the northern and southern stock subunits. 2024 temperature in the northern subunit (north of Hudson Canyon) was colder than black sea bass''s lower threshold of 8C.'
$ws.Range('G10').Value = 'du Pontavice, H., Miller, T. J., Stock, B. C., Chen, Z., & Saba, V. S. (2022). Ocean model-based covariates improve a marine fish stock assessment when observations are limited. ICES Journal of Marine Science, 79(4), 1259?1273. Jean-Michel, L., Eric, G., Romain, B.-B., Gilles, G., Ang?lique, M., Marie, D., Cl?ment, B., Mathieu, H., Olivier, L. G., Charly, R., Tony, C., Charles-Emmanuel, T., Florent, G., Giovanni, R., Mounir, B., Yann, D., & Pierre-Yves, L. T. (2021). The Copernicus Global 1/12? Oceanic and Sea Ice GLORYS12 Reanalysis. Frontiers in Earth Science, 9, 698876. https://doi.org/10.3389/feart.2021.698876'
$ws.Range('H10').Value = 'SOUTH'

# Move the active selection to match the saved view
$ws.Range('A10').Select()
